# Normalize the sheet/workbook to match the "committed" Excel-authored
# version of this fixture:
#   - rename the worksheet from "Sheet 1" to "Sheet1"
#   - bold + center the header row (A1:F1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
